$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 82 - new entry: multiband compressor intro
$ws.Cells(82,1).Value = "ok I made a multiband compressor. Before adding it,"
$ws.Cells(82,2).Value = 20544

# Row 83 - with MBC
$ws.Cells(83,1).Value = "with MBC"
$ws.Cells(83,2).Value = 20892
$ws.Cells(83,6).Value = "wow not bad tbh"

# Row 84 - MBC but not Comp
$ws.Cells(84,1).Value = "MBC but not Comp"
$ws.Cells(84,2).Value = 20720

# Row 85 - with neither
$ws.Cells(85,1).Value = "with neither"
$ws.Cells(85,2).Value = 20088

# Row 86 - with MBC (again)
$ws.Cells(86,1).Value = "with MBC"
$ws.Cells(86,2).Value = 20892
$ws.Cells(86,6).Value = "so the multiband comp adds 800 bytes to final size"

# Row 87 - with neither (again)
$ws.Cells(87,1).Value = "with neither"
$ws.Cells(87,2).Value = 20088

# Row 88 - Comp but not MBC
$ws.Cells(88,1).Value = "Comp but not MBC"
$ws.Cells(88,2).Value = 20544
$ws.Cells(88,6).Value = "normal comp is 450 bytes. So MB does add mor than I expected. I'm going to keep it though."

# Row 89 - MBC but not Comp (again)
$ws.Cells(89,1).Value = "MBC but not Comp"
$ws.Cells(89,2).Value = 20720

# Row 90 - optimizations in sat
$ws.Cells(90,1).Value = "optimizations in sat"
$ws.Cells(90,2).Value = 20664

# Row 91 - readd biquad
$ws.Cells(91,1).Value = "and considering this is supposed to *replace* comp, readd biquad"
$ws.Cells(91,2).Value = 20688
$ws.Cells(91,6).Value = "also optimizations in sat"

# Row 92 - without sat
$ws.Cells(92,1).Value = "without sat"
$ws.Cells(92,2).Value = 20204
$ws.Cells(92,6).Value = "so sat is now 480 bytes of code, that's pretty good tbh"

# Row 93 - closing remark (no B value set on purpose)
$ws.Cells(93,1).Value = "something to keep in mind is that sat could be built into MBC"

# Update selection to match the author's final cursor position
$ws.Range("A94").Select()
